$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (Veneto) at row 23
$ws.Range("A23").Value = "Veneto"
$ws.Range("B23").Value = 2019
$ws.Range("C23").Value = 56437283
$ws.Range("D23").Value = 29907
$ws.Range("E23").Value = 69026517
$ws.Range("F23").Value = 294762
$ws.Range("G23").Value = 429550
$ws.Range("H23").Value = 765245

# Update selection to match post-edit state
$ws.Range("A23:H23").Select()
